# Commit: Change sheet name of files
# - Rename the "Recordings" sheet to "Files"
# - Make the "Files" sheet the active/selected tab (was "Collection")
# - Adjust row heights for header rows 1 and 2 on the "Files" sheet

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Recordings")
$ws.Name = "Files"

# Row height tweaks (points)
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 17.25

# Make "Files" the active sheet/tab (was "Collection")
$ws.Activate()
